$wb = $excel.ActiveWorkbook

# Duplicate the "Croatia" sheet (same template as every market tab) and drop
# the copy right after it, at the end of the tab strip.
$croatia = $wb.Worksheets.Item("Croatia")
[void]$croatia.Copy($null, $croatia)

# Copying moves the selection/active-tab onto the new sheet; put Croatia's
# view back the way it ends up after the focus has moved on (whole sheet
# selected, no longer the active tab).
[void]$croatia.Activate()
[void]$croatia.Cells.Select()

# The copy is inserted immediately after Croatia, i.e. it's now the last tab.
$greece = $wb.Worksheets.Item($wb.Worksheets.Count)
$greece.Name = "Greece"
$greece.Range("B4").Value = "NGC-4119/T3163"
$greece.Range("B2").Value = "Greece Market"

# Leave Greece as the active sheet/selection.
[void]$greece.Activate()
[void]$greece.Range("B2").Select()
